$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04339299999999999
$ws.Range("H2").Value = 0.130179
$ws.Range("I2").Value = 0.0698021577815419
$ws.Range("J2").Value = 0.0698021577815419
$ws.Range("M2").Value = 35.991783
$ws.Range("N2").Value = 107.975349
$ws.Range("O2").Value = 0.3909505149237033
$ws.Range("P2").Value = 0.3909505149237033
$ws.Range("Q2").Value = 1.561791439719
$ws.Range("R2").Value = 14.056122957471
$ws.Range("S2").Value = 0.02728918952747939
$ws.Range("T2").Value = 0.02728918952747939

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04339299999999999
$ws.Range("H3").Value = 0.130179
$ws.Range("I3").Value = 0.0698021577815419
$ws.Range("J3").Value = 0.0698021577815419
$ws.Range("O3").Value = 0.537552751174421
$ws.Range("P3").Value = 0.537552751174421
$ws.Range("Q3").Value = 2.147446423866333
$ws.Range("R3").Value = 19.327017814797
$ws.Range("S3").Value = 0.03752234195337887
$ws.Range("T3").Value = 0.03752234195337887

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04339299999999999
$ws.Range("H4").Value = 0.130179
$ws.Range("I4").Value = 0.0698021577815419
$ws.Range("J4").Value = 0.0698021577815419
$ws.Range("M4").Value = 6.58215
$ws.Range("O4").Value = 0.07149673390187571
$ws.Range("P4").Value = 0.07149673390187571
$ws.Range("Q4").Value = 0.2856192349499999
$ws.Range("R4").Value = 2.57057311455
$ws.Range("S4").Value = 0.004990626300683644
$ws.Range("T4").Value = 0.004990626300683644

# Row 5
$ws.Range("I5").Value = 0.6764796878879081
$ws.Range("J5").Value = 0.6764796878879081
$ws.Range("M5").Value = 35.991783
$ws.Range("N5").Value = 107.975349
$ws.Range("O5").Value = 0.3909505149237033
$ws.Range("P5").Value = 0.3909505149237033
$ws.Range("Q5").Value = 15.135924436515
$ws.Range("R5").Value = 136.223319928635
$ws.Range("S5").Value = 0.2644700823152038
$ws.Range("T5").Value = 0.2644700823152038

# Row 6
$ws.Range("I6").Value = 0.6764796878879081
$ws.Range("J6").Value = 0.6764796878879081
$ws.Range("O6").Value = 0.537552751174421
$ws.Range("P6").Value = 0.537552751174421
$ws.Range("S6").Value = 0.3636435173377587
$ws.Range("T6").Value = 0.3636435173377587

# Row 7
$ws.Range("I7").Value = 0.6764796878879081
$ws.Range("J7").Value = 0.6764796878879081
$ws.Range("M7").Value = 6.58215
$ws.Range("O7").Value = 0.07149673390187571
$ws.Range("P7").Value = 0.07149673390187571
$ws.Range("Q7").Value = 2.768046390749999
$ws.Range("S7").Value = 0.0483660882349457
$ws.Range("T7").Value = 0.0483660882349457

# Row 8
$ws.Range("I8").Value = 0.2537181543305499
$ws.Range("J8").Value = 0.2537181543305499
$ws.Range("M8").Value = 35.991783
$ws.Range("N8").Value = 107.975349
$ws.Range("O8").Value = 0.3909505149237033
$ws.Range("P8").Value = 0.3909505149237033
$ws.Range("Q8").Value = 5.676827968197
$ws.Range("R8").Value = 51.091451713773
$ws.Range("S8").Value = 0.09919124308102013
$ws.Range("T8").Value = 0.09919124308102013

# Row 9
$ws.Range("I9").Value = 0.2537181543305499
$ws.Range("J9").Value = 0.2537181543305499
$ws.Range("O9").Value = 0.537552751174421
$ws.Range("P9").Value = 0.537552751174421
$ws.Range("S9").Value = 0.1363868918832835
$ws.Range("T9").Value = 0.1363868918832835

# Row 10
$ws.Range("I10").Value = 0.2537181543305499
$ws.Range("J10").Value = 0.2537181543305499
$ws.Range("M10").Value = 6.58215
$ws.Range("R10").Value = 9.343565971649999
$ws.Range("S10").Value = 0.01814001936624636
$ws.Range("T10").Value = 0.01814001936624636

